# NIT-9017819120.xlsx — "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The worker detail table (rows 16-39) is re-sorted/grouped by worker
# (Rolando Perez Herrera, Deivi Jose Torres Pedroza, Coby Brayan Maury
# Rodriguez), each with periods 2412 down to 2406 (descending), and the
# prior 2501 period rows are dropped entirely (8 periods -> 7 periods,
# 24 data rows -> 21 data rows). The header summary (VALOR MORA, Cant.
# Periodos) is refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- refresh header summary values -----------------------------------
$ws.Range("E11").Value = 1064266
$ws.Range("F13").Value = 7

# --- rewrite the worker detail table (rows 16-36) ---------------------
# Columns: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo Mora, F=Valor Mora, G=Salario Basico
$data = @(
    @("CC","7920637","ROLANDO PEREZ HERRERA","2412",52000,1300000),
    @("CC","7920637","ROLANDO PEREZ HERRERA","2411",52000,1300000),
    @("CC","7920637","ROLANDO PEREZ HERRERA","2410",52000,1300000),
    @("CC","7920637","ROLANDO PEREZ HERRERA","2409",52000,1300000),
    @("CC","7920637","ROLANDO PEREZ HERRERA","2408",52000,1300000),
    @("CC","7920637","ROLANDO PEREZ HERRERA","2407",52000,1300000),
    @("CC","7920637","ROLANDO PEREZ HERRERA","2406",46800,1300000),
    @("CC","1047389355","DEIVI JOSE TORRES PEDROZA","2412",52000,1300000),
    @("CC","1047389355","DEIVI JOSE TORRES PEDROZA","2411",52000,1300000),
    @("CC","1047389355","DEIVI JOSE TORRES PEDROZA","2410",52000,1300000),
    @("CC","1047389355","DEIVI JOSE TORRES PEDROZA","2409",52000,1300000),
    @("CC","1047389355","DEIVI JOSE TORRES PEDROZA","2408",52000,1300000),
    @("CC","1047389355","DEIVI JOSE TORRES PEDROZA","2407",52000,1300000),
    @("CC","1047389355","DEIVI JOSE TORRES PEDROZA","2406",46800,1300000),
    @("CC","1007120657","COBY BRAYAN MAURY RODRIGUEZ","2412",52000,1300000),
    @("CC","1007120657","COBY BRAYAN MAURY RODRIGUEZ","2411",52000,1300000),
    @("CC","1007120657","COBY BRAYAN MAURY RODRIGUEZ","2410",52000,1300000),
    @("CC","1007120657","COBY BRAYAN MAURY RODRIGUEZ","2409",52000,1300000),
    @("CC","1007120657","COBY BRAYAN MAURY RODRIGUEZ","2408",52000,1300000),
    @("CC","1007120657","COBY BRAYAN MAURY RODRIGUEZ","2407",52000,1300000),
    @("CC","1007120657","COBY BRAYAN MAURY RODRIGUEZ","2406",34666,1300000)
)

$r = 16
foreach ($row in $data) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}

# --- drop the now-unused trailing rows (old rows 37-39 held period 2501) --
# Deleting shifts the footer (signature block) rows up from 44/45 to 41/42.
$ws.Rows("37:39").Delete()
